$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "66.902.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value2 = "  -3.80%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "3.338.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value2 = "  -0.96%  "
$ws.Range("E4").Value2 = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "574.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "  -3.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "183.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = "  -4.54%  "
$ws.Range("E7").Value2 = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.598"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value2 = "  -1.98%  "
$ws.Range("E9").Value2 = "  -3.79%  "
$ws.Range("E10").Value2 = "  -1.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "0.404"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value2 = "  -4.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "3.918.36"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value2 = "  -1.02%  "
$ws.Range("E13").Value2 = "  -0.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "27.19"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value2 = "  -5.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "66.906.58"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value2 = "  -3.81%  "
$ws.Range("E16").Value2 = "  -2.63%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "3.334.38"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value2 = "  -1.10%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "435.84"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value2 = "  -2.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "13.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = "  -1.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "5.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = "  -2.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "7.66"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value2 = "  -2.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "73.83"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value2 = "  +0.34%  "
$ws.Range("E23").Value2 = "  +0.04%  "
$ws.Range("E24").Value2 = "  +0.24%  "
$ws.Range("E25").Value2 = "  -2.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "0.191"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value2 = "  -0.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "9.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value2 = "  -5.42%  "
$ws.Range("E28").Value2 = "  -0.05%  "
$ws.Range("E29").Value2 = "  -1.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "22.85"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value2 = "  -1.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "5.34"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value2 = "  -5.12%  "
$ws.Range("B32").Value2 = "USDe"
$ws.Range("C32").Value2 = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "0.999"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value2 = "  +0.02%  "
$ws.Range("B33").Value2 = "Aptos"
$ws.Range("C33").Value2 = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "6.85"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value2 = "  -2.94%  "
$ws.Range("E34").Value2 = "  -4.74%  "
$ws.Range("E35").Value2 = "  -0.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "160.16"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value2 = "  -2.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "27.42"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value2 = "  +0.34%  "
$ws.Range("E38").Value2 = "  -4.86%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "2.842.73"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value2 = "  +3.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "0.793"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value2 = "  -3.60%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "4.45"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value2 = "  -3.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "6.25"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value2 = "  -4.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "0.0676"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value2 = "  -2.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "40.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value2 = "  -1.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "24.65"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value2 = "  -3.92%  "
$ws.Range("E46").Value2 = "  -6.81%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "322.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value2 = "  -6.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "0.0273"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value2 = "  -4.70%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "0.989"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value2 = "  -4.50%  "
$ws.Range("B50").Value2 = "Arweave"
$ws.Range("C50").Value2 = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "30.98"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value2 = "  -5.81%  "
$ws.Range("B51").Value2 = "Cosmos"
$ws.Range("C51").Value2 = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "6.16"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value2 = "  -2.91%  "
